# Weekly update: a new week of "Ciboulette" price observations is inserted
# at the top of the data block (rows 178-179), pushing every existing
# observation down by two rows (old row 178 -> new row 180, ... old row 254
# -> new row 256). The used range grows from A1:R254 to A1:R256 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 178-179; everything below (178..254) shifts
# down to 180..256, keeping its data intact (Excel copies the formatting of
# the row above, which matches the existing date-style column D already
# uses throughout the sheet).
$ws.Range("A178:A179").EntireRow.Insert()

# Row 178 - "Primera" quality observation for the new week (Fecha 44466)
$ws.Cells.Item(178, 1).Value = 9
$ws.Cells.Item(178, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(178, 3).Value = "Metropolitana"
$ws.Cells.Item(178, 4).Value = 44466
$ws.Cells.Item(178, 5).Value = 13
$ws.Cells.Item(178, 6).Value = 100112039
$ws.Cells.Item(178, 7).Value = "Ciboulette"
$ws.Cells.Item(178, 8).Value = "Sin especificar"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 250
$ws.Cells.Item(178, 11).Value = 1000
$ws.Cells.Item(178, 12).Value = 1300
$ws.Cells.Item(178, 13).Value = 1150
$ws.Cells.Item(178, 14).Value = "$/docena de atados"
$ws.Cells.Item(178, 15).Value = "Región Metropolitana"
$ws.Cells.Item(178, 16).Value = 383
$ws.Cells.Item(178, 17).Value = 3
$ws.Cells.Item(178, 18).Value = "Hortaliza"

# Row 179 - "Segunda" quality observation for the same new week (Fecha 44466)
$ws.Cells.Item(179, 1).Value = 9
$ws.Cells.Item(179, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(179, 3).Value = "Metropolitana"
$ws.Cells.Item(179, 4).Value = 44466
$ws.Cells.Item(179, 5).Value = 13
$ws.Cells.Item(179, 6).Value = 100112039
$ws.Cells.Item(179, 7).Value = "Ciboulette"
$ws.Cells.Item(179, 8).Value = "Sin especificar"
$ws.Cells.Item(179, 9).Value = "Segunda"
$ws.Cells.Item(179, 10).Value = 106
$ws.Cells.Item(179, 11).Value = 900
$ws.Cells.Item(179, 12).Value = 900
$ws.Cells.Item(179, 13).Value = 900
$ws.Cells.Item(179, 14).Value = "$/docena de atados"
$ws.Cells.Item(179, 15).Value = "Región Metropolitana"
$ws.Cells.Item(179, 16).Value = 300
$ws.Cells.Item(179, 17).Value = 3
$ws.Cells.Item(179, 18).Value = "Hortaliza"
